$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F6").Value = 419
$ws1.Range("F8").Value = 13305
$ws1.Range("F10").Value = 47
$ws1.Range("F11").Value = 5388
$ws1.Range("F13").Value = 29
$ws1.Range("F16").Value = 1210
$ws1.Range("F18").Value = 146
$ws1.Range("F21").Value = 7262
$ws1.Range("F23").Value = 3663

$ws4.Range("F7").Value = 419
$ws4.Range("F9").Value = 13305
$ws4.Range("F11").Value = 47
$ws4.Range("F12").Value = 5388
$ws4.Range("F14").Value = 29
$ws4.Range("F17").Value = 1210
$ws4.Range("F19").Value = 146
$ws4.Range("F23").Value = 7262
$ws4.Range("F25").Value = 3663
